$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right below the header row so the existing data
# (originally rows 2-18) shifts down to rows 4-20.
$ws.Range("A2:A3").EntireRow.Insert()

# New row 2: 하이젠알앤엠  (leading ' forces text so date-like strings are not
# auto-converted into Excel date serials - matches the rest of the sheet,
# where every column is plain text/number with no special formatting)
$row2 = @("'2024-06-18", "하이젠알앤엠", "한국", "'2024-06-21", "'2024-06-27", 23800000, 3400000, "-", 4500, 5500, "-", 7000, "-", "-", 0, "-", "-", "2549.13 : 1", "-", "-")
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

# New row 3: 한국제15호스팩
$row3 = @("'2024-06-17", "한국제15호스팩", "한국", "'2024-06-20", "'2024-06-26", 12500000, 6250000, "-", 2000, 2000, "-", 2000, "-", "-", 0, "-", "-", "736.67 : 1", "-", "-")
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

# The row-insert copies formatting (bold/border/center) from the header row
# above; strip it so the new rows plainly match the rest of the data rows.
$ws.Range("A2:T3").ClearFormats()
